# Revert "Merging 0.1.8 w VitalSigns"
#
# - Rename the "Include #0" sheet to "Include from SNOMED CT".
# - On the "Metadata" sheet, roll several published-metadata values back
#   to their pre-merge values (Version, Status, Date, Contact) and drop
#   the "Jurisdiction" row entirely (rows below it shift up).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value  = "0.1.6"
$ws1.Range("B6").Value  = "active"
$ws1.Range("B8").Value  = "2023-05-05T10:50:04-05:00"
$ws1.Range("B10").Value = "No display for ContactDetail"
$ws1.Range("B11").Value = "No display for ContactDetail"

# Row 12 was "Jurisdiction" | (blank) - remove it, shifting later rows up.
$ws1.Rows.Item(12).Delete()

$ws2 = $wb.Worksheets.Item("Include #0")
$ws2.Name = "Include from SNOMED CT"
